$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 18 de Abril de 2020 a las 23:52"

# Germany (row 8) - updated case counts
$ws.Range("B8").Value = 143342
$ws.Range("C8").Value = 1945
$ws.Range("E8").Value = 53483
$ws.Range("G8").Value = 107
$ws.Range("H8").Value = 4459

# Niger (row 95) - updated case counts
$ws.Range("B95").Value = 639
$ws.Range("C95").Value = 12
$ws.Range("D95").Value = 113
$ws.Range("E95").Value = 507
$ws.Range("G95").Value = 1
$ws.Range("H95").Value = 19

# Row 128 was Tanzania, now becomes Guadalupe (with updated counts, overtaking Tanzania in ranking)
$ws.Range("A128").Value = "Guadalupe"
$ws.Range("B128").Value = 148
$ws.Range("C128").Value = 3
$ws.Range("D128").Value = 73
$ws.Range("E128").Value = 67
$ws.Range("F128").Value = 13
$ws.Range("H128").Value = 8

# Row 129 was Guadalupe, now becomes Tanzania (keeps the old Tanzania counts)
$ws.Range("A129").Value = "Tanzania"
$ws.Range("B129").Value = 147
$ws.Range("D129").Value = 11
$ws.Range("E129").Value = 131
$ws.Range("F129").Value = 4
$ws.Range("H129").Value = 5
